$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.384.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "1.721.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'242.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.9992"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4883"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "1.727.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.06981"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'15.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.517"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.5972"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'77.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "26.387.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.9990"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.000007156"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "1.946.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.441"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'8.486"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'5.081"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'138.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'1.403"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'106.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.728"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.91%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'3.905"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.08032"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'3.652"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.6216"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.9187"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.963"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.381"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.9985"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.01476"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'99.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'5.419"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.3836"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'6.889"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.1161"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.05365"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'30.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'7.666"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'51.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.216"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.65%  "
$ws.Range("E51").Style = "Normal"
